# Cards workbook edit: add a "Nom" column (card position 1-9 within its
# suit block) to the Feuil1 table.
#
# The header row gets a clean column insert before column B (old B/C/D/E
# "Valeur"/"PointsSansAtout"/"PointsAvecAtout"/"Image" headers shift one
# column right, and the new B1 becomes "Nom").
#
# For the data rows, the old "Valeur" column (B) stays put, and the new
# "Nom" numbering (1..9, repeating every 9 rows - one per suit) is written
# into column C, with the old C/D/E values (PointsSansAtout/PointsAvecAtout
# /Image) shifting one column right to D/E/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B. This shifts the whole sheet's
# B:E columns to C:F for every row (dimension grows to F37 automatically).
[void]$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Nom"

# Fix up the data rows (2-37): the insert pushed the old "Valeur" value for
# each row into column C; move it back to B, then write the sequential
# "Nom" position (1-9 per suit block of 9 rows) into C.
for ($r = 2; $r -le 37; $r++) {
    $valeur = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $valeur
    $nom = (($r - 2) % 9) + 1
    $ws.Cells.Item($r, 3).Value = $nom
}

# Match the author's final selection/active cell.
[void]$ws.Range("I8").Select()
